{"js": "const sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst section = sections.items[0];\nconst header = section.getHeader(\"Primary\");\nconst footer = section.getFooter(\"Primary\");\nheader.clear();\nfooter.clear();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Add a default header and footer to the document's (first/primary) section.\n# Setting the paragraph Style on the header/footer Range materializes the\n# header1.xml / footer1.xml parts (with a single empty paragraph styled\n# \"Header\" / \"Footer\") and wires up the <w:headerReference>/<w:footerReference>\n# in the section's sectPr, without forcing creation of first-page/even-page\n# variants.\n$section = $d.Sections.Item(1)\n\n$header = $section.Headers.Item(1)\n$header.Range.Paragraphs.Item(1).Style = \"Header\"\n\n$footer = $section.Footers.Item(1)\n$footer.Range.Paragraphs.Item(1).Style = \"Footer\"\n"}
